$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3, and bump the "Förändrad" date by 1 day
$ws.Range("A2").Value = "A 36523-2022"
$ws.Range("C2").Value = 46079
$ws.Range("G2").Value = 0.2

$ws.Range("A3").Value = "A 36578-2022"
$ws.Range("C3").Value = 46079
$ws.Range("G3").Value = 0.3
